$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scale the "push and pull" fund data in column D (rows 2-33) by 10000
$lastRow = 33
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 4)
    $current = $cell.Value()
    if ($current -ne $null) {
        $cell.Value = $current * 10000
    }
}
